$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.982.34"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "3.214.88"
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "603.89"
$ws.Range("E5").Value = "  +4.30%  "

$ws.Range("D6").Value = "153.89"
$ws.Range("E6").Value = "  +2.24%  "

$ws.Range("D8").Value = "3.214.91"
$ws.Range("E8").Value = "  +1.60%  "

$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  -1.22%  "

$ws.Range("D12").Value = "0.511"
$ws.Range("E12").Value = "  +2.16%  "

$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").Value = "39.49"
$ws.Range("E14").Value = "  +5.39%  "

$ws.Range("D15").Value = "3.743.78"
$ws.Range("E15").Value = "  +1.71%  "

$ws.Range("D16").Value = "7.52"
$ws.Range("E16").Value = "  +5.15%  "

$ws.Range("D17").Value = "66.124.70"
$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Value = "3.219.41"
$ws.Range("E18").Value = "  +1.54%  "

$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").Value = "511.22"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("D21").Value = "15.47"
$ws.Range("E21").Value = "  +4.24%  "

$ws.Range("D22").Value = "0.739"
$ws.Range("E22").Value = "  +1.88%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "15.49"
$ws.Range("E23").Value = "  +1.18%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").Value = "'8.10"
$ws.Range("E24").Value = "  +4.11%  "

$ws.Range("D25").Value = "85.14"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  +2.32%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "3.02"
$ws.Range("E28").Value = "  +3.44%  "

$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  +4.79%  "

$ws.Range("D30").Value = "2.86"
$ws.Range("E30").Value = "  +2.95%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "28.21"
$ws.Range("E31").Value = "  +1.72%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.80"
$ws.Range("E32").Value = "  +7.94%  "

$ws.Range("E33").Value = "  +2.34%  "

$ws.Range("E34").Value = "  +0.15%  "

$ws.Range("D35").Value = "6.61"
$ws.Range("E35").Value = "  +0.55%  "

$ws.Range("D36").Value = "55.19"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").Value = "0.0906"
$ws.Range("E37").Value = "  +1.53%  "

$ws.Range("D38").Value = "485.09"
$ws.Range("E38").Value = "  +2.82%  "

$ws.Range("D39").Value = "'0.0420"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").Value = "2.95"
$ws.Range("E40").Value = "  -5.37%  "

$ws.Range("D41").Value = "8.94"
$ws.Range("E41").Value = "  +3.60%  "

$ws.Range("D42").Value = "0.301"
$ws.Range("E42").Value = "  +6.39%  "

$ws.Range("E43").Value = "  +1.82%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.946.75"
$ws.Range("E44").Value = "  -3.75%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.45"
$ws.Range("E45").Value = "  +1.50%  "

$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0₃0643"
$ws.Range("E46").Value = "  +6.57%  "

$ws.Range("D47").Value = "28.68"
$ws.Range("E47").Value = "  -2.12%  "

$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  +3.31%  "

$ws.Range("D51").Value = "120.73"
$ws.Range("E51").Value = "  +0.85%  "
